# houstonNumbers.xlsx update — "Updated excel with data from 3/23/2020"
#
# 1) Two existing "age bracket" cells (C16, C130) were mis-typed back when the
#    sheet was first built: they hold the date serial 44124 (formatted d-mmm,
#    which *looked* like "20-Oct") where the intended value was the text
#    bracket "10-20". Fix those to the text "10-20" and make their number
#    format Text so it can't be re-interpreted as a date again.
# 2) Append the 3/22 and 3/23/2020 case rows (144-190).
# 3) Cosmetic: column widths for the now-visible D/E columns, and leave the
#    selection where the author left off (E181, scrolled to show row 148+).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix C16 / C130: numeric date-serial -> text "10-20"
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "10-20"
$ws.Range("C130").Value = "10-20"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C130").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 2) Append rows 144-190
# ---------------------------------------------------------------------------

# Column A holds dates; copy the number format from the last existing date
# cell (A143) onto the new range so the new cells share the real date style
# instead of picking up General formatting.
$ws.Range("A143").Copy()
$ws.Range("A144:A190").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A144:A157").Value = 43912
$ws.Range("A158:A190").Value = 43913

$ws.Range("B144").Value = "F"
$ws.Range("C144").Value = "40-50"
$ws.Range("D144").Value = "Galveston"
$ws.Range("E144").Value = "Travel"
$ws.Range("B145").Value = "F"
$ws.Range("C145").Value = "30-40"
$ws.Range("D145").Value = "Harris"
$ws.Range("E145").Value = "Community Spread"
$ws.Range("B146").Value = "M"
$ws.Range("C146").Value = "50-60"
$ws.Range("D146").Value = "Harris"
$ws.Range("E146").Value = "Exposed"
$ws.Range("B147").Value = "M"
$ws.Range("C147").Value = "40-50"
$ws.Range("D147").Value = "Harris"
$ws.Range("E147").Value = "Community Spread"
$ws.Range("B148").Value = "M"
$ws.Range("C148").Value = "20-30"
$ws.Range("D148").Value = "Harris"
$ws.Range("E148").Value = "Community Spread"
$ws.Range("B149").Value = "M"
$ws.Range("C149").Value = "40-50"
$ws.Range("D149").Value = "Harris"
$ws.Range("E149").Value = "Community Spread"
$ws.Range("B150").Value = "M"
$ws.Range("C150").Value = "30-40"
$ws.Range("D150").Value = "Houston"
$ws.Range("E150").Value = "Travel"
$ws.Range("B151").Value = "F"
$ws.Range("C151").Value = "40-50"
$ws.Range("D151").Value = "Montgomery"
$ws.Range("E151").Value = "Travel"
$ws.Range("B152").Value = "M"
$ws.Range("C152").Value = "40-50"
$ws.Range("D152").Value = "Montgomery"
$ws.Range("E152").Value = "Community Spread"
$ws.Range("B153").Value = "M"
$ws.Range("C153").Value = "70-80"
$ws.Range("D153").Value = "Galveston"
$ws.Range("E153").Value = "Travel"
$ws.Range("B154").Value = "F"
$ws.Range("C154").Value = "20-30"
$ws.Range("D154").Value = "Galveston"
$ws.Range("E154").Value = "Community Spread"
$ws.Range("B155").Value = "F"
$ws.Range("C155").Value = "50-60"
$ws.Range("D155").Value = "Galveston"
$ws.Range("E155").Value = "Community Spread"
$ws.Range("D156").Value = "Brazos"
$ws.Range("E156").Value = "Community Spread"
$ws.Range("D157").Value = "Brazos"
$ws.Range("E157").Value = "Community Spread"
$ws.Range("D158").Value = "Fort Bend"
$ws.Range("D159").Value = "Fort Bend"
$ws.Range("D160").Value = "Fort Bend"
$ws.Range("D161").Value = "Fort Bend"
$ws.Range("D162").Value = "Fort Bend"
$ws.Range("D163").Value = "Fort Bend"
$ws.Range("D164").Value = "Fort Bend"
$ws.Range("D165").Value = "Fort Bend"
$ws.Range("D166").Value = "Fort Bend"
$ws.Range("D167").Value = "Fort Bend"
$ws.Range("D168").Value = "Fort Bend"
$ws.Range("D169").Value = "Fort Bend"
$ws.Range("D170").Value = "Fort Bend"
$ws.Range("B171").Value = "M"
$ws.Range("C171").Value = "60-70"
$ws.Range("D171").Value = "Brazoria"
$ws.Range("E171").Value = "Community Spread"
$ws.Range("B172").Value = "M"
$ws.Range("C172").Value = "30-40"
$ws.Range("D172").Value = "Brazoria"
$ws.Range("B173").Value = "F"
$ws.Range("C173").Value = "0-10"
$ws.Range("D173").Value = "Galveston"
$ws.Range("E173").Value = "Travel"

# C174 is the third occurrence of the mis-typed "10-20" date-serial bracket —
# same fix as C16 / C130 above.
$ws.Range("B174").Value = "F"
$ws.Range("C174").Value = "10-20"
$ws.Range("C174").NumberFormat = "@"
$ws.Range("D174").Value = "Montgomery"
$ws.Range("E174").Value = "Travel"

$ws.Range("B175").Value = "F"
$ws.Range("C175").Value = "20-30"
$ws.Range("D175").Value = "Montgomery"
$ws.Range("E175").Value = "Community Spread"
$ws.Range("B176").Value = "M"
$ws.Range("C176").Value = "50-60"
$ws.Range("D176").Value = "Montgomery"
$ws.Range("E176").Value = "Travel"
$ws.Range("B177").Value = "M"
$ws.Range("C177").Value = "50-60"
$ws.Range("D177").Value = "Montgomery"
$ws.Range("E177").Value = "Travel"
$ws.Range("B178").Value = "F"
$ws.Range("C178").Value = "40-50"
$ws.Range("D178").Value = "Harris"
$ws.Range("B179").Value = "F"
$ws.Range("C179").Value = "60-70"
$ws.Range("D179").Value = "Harris"
$ws.Range("E179").Value = "Community Spread"
$ws.Range("B180").Value = "M"
$ws.Range("C180").Value = "50-60"
$ws.Range("D180").Value = "Harris"
$ws.Range("E180").Value = "Travel"
$ws.Range("B181").Value = "F"
$ws.Range("C181").Value = "20-30"
$ws.Range("D181").Value = "Houston"
$ws.Range("E181").Value = "Community Spread"
$ws.Range("B182").Value = "M"
$ws.Range("C182").Value = "50-60"
$ws.Range("D182").Value = "Matagorda"
$ws.Range("B183").Value = "M"
$ws.Range("C183").Value = "20-30"
$ws.Range("D183").Value = "Matagorda"
$ws.Range("B184").Value = "M"
$ws.Range("C184").Value = "20-30"
$ws.Range("D184").Value = "Galveston"
$ws.Range("E184").Value = "Travel"
$ws.Range("B185").Value = "M"
$ws.Range("C185").Value = "20-30"
$ws.Range("D185").Value = "Galveston"
$ws.Range("E185").Value = "Travel"
$ws.Range("B186").Value = "M"
$ws.Range("C186").Value = "40-50"
$ws.Range("D186").Value = "Galveston"
$ws.Range("E186").Value = "Travel"
$ws.Range("B187").Value = "M"
$ws.Range("C187").Value = "50-60"
$ws.Range("D187").Value = "Galveston"
$ws.Range("E187").Value = "Travel"
$ws.Range("B188").Value = "M"
$ws.Range("C188").Value = "60-70"
$ws.Range("D188").Value = "Galveston"
$ws.Range("E188").Value = "Travel"
$ws.Range("B189").Value = "F"
$ws.Range("C189").Value = "20-30"
$ws.Range("D189").Value = "Brazoria"
$ws.Range("E189").Value = "Community Spread"
$ws.Range("F189").Value = "Y"
$ws.Range("B190").Value = "M"
$ws.Range("C190").Value = "50-60"
$ws.Range("D190").Value = "Brazoria"
$ws.Range("E190").Value = "Community Spread"

# ---------------------------------------------------------------------------
# 3) Cosmetics: column widths + leave the selection/scroll where the author
#    ended up after entering the new data.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 10.83
$ws.Columns.Item(5).ColumnWidth = 16

$ws.Range("E181").Select()
